# Replace every addition/subtraction problem in the single 20x5 table with
# the updated set of problems (diff from commit "Update master to output
# generated at c8c62b6"). Cells are addressed directly by (row, col) via
# Table.Cell(r, c).Range.Text rather than Find/Replace, because several of
# the original expressions (e.g. "37+11=") are duplicated in the table and
# a text-based Find/Replace could not disambiguate which occurrence to
# update.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    @("55+16=", "39+25=", "77-53=", "38-11=", "31+5="),
    @("35-7=", "93-83=", "16+71=", "90-81=", "81+0="),
    @("89-86=", "11+5=", "17+77=", "43+54=", "46+41="),
    @("31+46=", "0+66=", "86-24=", "36+4=", "32+39="),
    @("37+39=", "57+13=", "17+40=", "9+24=", "26+25="),
    @("45+4=", "55-16=", "35+44=", "97-16=", "35-0="),
    @("21+69=", "9-7=", "34-4=", "53-39=", "45-34="),
    @("23+22=", "52+25=", "97-22=", "59+25=", "39-26="),
    @("82-29=", "43-30=", "36+50=", "85-21=", "49+44="),
    @("30+55=", "75-20=", "55-11=", "84-12=", "15+49="),
    @("64-22=", "74-10=", "2+88=", "93-71=", "13-6="),
    @("34+50=", "59+31=", "65-41=", "19+59=", "88-60="),
    @("3+34=", "0+96=", "54+20=", "85-28=", "20+11="),
    @("10+43=", "48-6=", "7+15=", "86-72=", "58-12="),
    @("54-33=", "1+49=", "51-4=", "25-20=", "59-19="),
    @("34+46=", "6+47=", "40+44=", "63-25=", "88-73="),
    @("14-4=", "15+75=", "44-29=", "41-3=", "65-25="),
    @("95-3=", "14+43=", "47-21=", "2+56=", "40+25="),
    @("71-62=", "90-88=", "52-46=", "62+27=", "76-10="),
    @("58-11=", "84-12=", "10+69=", "28-3=", "43+36=")
)

if ($t.Rows.Count -ne $values.Length -or $t.Columns.Count -ne $values[0].Length) {
    Write-Host "Warning: table shape" $t.Rows.Count "x" $t.Columns.Count "does not match expected" $values.Length "x" $values[0].Length
}

for ($r = 1; $r -le $values.Length; $r++) {
    $rowVals = $values[$r - 1]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $t.Cell($r, $c).Range.Text = $rowVals[$c - 1]
    }
}

Write-Host "Done updating table cells."